# Updates cryptos list values to match the latest scrape (GitHub Actions run).
# Cells in column D that are numeric-looking strings are forced to Text format
# first so Excel keeps exact precision (trailing zeros, very small decimals, etc.)
# instead of silently converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.386.60"
$ws.Range("E2").Value = "  -4.09%  "
$ws.Range("D3").Value = "2.980.11"
$ws.Range("E3").Value = "  -5.90%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "560.92"
$ws.Range("E5").Value = "  -5.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "125.16"
$ws.Range("E6").Value = "  -7.08%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "2.976.47"
$ws.Range("E8").Value = "  -5.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.494"
$ws.Range("E9").Value = "  -3.92%  "
$ws.Range("E10").Value = "  -5.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.12"
$ws.Range("E11").Value = "  -2.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.435"
$ws.Range("E12").Value = "  -4.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000223"
$ws.Range("E13").Value = "  -6.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.46"
$ws.Range("E14").Value = "  -6.32%  "
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("D16").Value = "3.470.87"
$ws.Range("E16").Value = "  -5.90%  "
$ws.Range("D17").Value = "60.563.91"
$ws.Range("E17").Value = "  -3.73%  "
$ws.Range("D18").Value = "2.988.42"
$ws.Range("E18").Value = "  -5.75%  "
$ws.Range("E19").Value = "  -6.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "428.83"
$ws.Range("E20").Value = "  -6.78%  "
$ws.Range("E21").Value = "  -6.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.656"
$ws.Range("E22").Value = "  -6.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.10"
$ws.Range("E23").Value = "  -6.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.87"
$ws.Range("E24").Value = "  -3.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "78.23"
$ws.Range("E25").Value = "  -6.36%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.47"
$ws.Range("E28").Value = "  -7.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.15"
$ws.Range("E29").Value = "  -7.02%  "
$ws.Range("E30").Value = "  -7.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.22"
$ws.Range("E31").Value = "  -7.15%  "
$ws.Range("E32").Value = "  -10.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0933"
$ws.Range("E33").Value = "  -9.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.26"
$ws.Range("E34").Value = "  -4.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.946"
$ws.Range("E35").Value = "  -8.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.50"
$ws.Range("E36").Value = "  -5.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "49.32"
$ws.Range("E37").Value = "  -3.63%  "
$ws.Range("D38").Value = "0.0₃0660"
$ws.Range("E38").Value = "  -6.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0357"
$ws.Range("E39").Value = "  -8.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.75"
$ws.Range("E40").Value = "  -4.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "374.43"
$ws.Range("E41").Value = "  -6.68%  "
$ws.Range("E42").Value = "  -4.93%  "
$ws.Range("D43").Value = "2.668.98"
$ws.Range("E43").Value = "  -4.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.43"
$ws.Range("E44").Value = "  -7.66%  "
$ws.Range("E46").Value = "  -6.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "119.00"
$ws.Range("E47").Value = "  -3.75%  "
$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.95"
$ws.Range("E48").Value = "  -7.42%  "
$ws.Range("B49").Value = "Arweave"
$ws.Range("C49").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "32.68"
$ws.Range("E49").Value = "  -4.70%  "
$ws.Range("E50").Value = "  -5.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.26"
$ws.Range("E51").Value = "  -8.11%  "

Write-Host "Applied cryptos list update (90 cell changes)"
